$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 61766.223
$ws.Range("I116").Value = 78334.71000000001
$ws.Range("J116").Value = 3776.5
$ws.Range("K116").Value = 78334.71000000001
$ws.Range("L116").Value = 3776.5
$ws.Range("M116").Value = -74892.71000000001
$ws.Range("N116").Value = -10660.5

$ws.Range("H138").Value = 4223.933
$ws.Range("I138").Value = 2251.1738
$ws.Range("J138").Value = 6286.364
$ws.Range("K138").Value = 6753.5214
$ws.Range("L138").Value = 18859.092
$ws.Range("M138").Value = -1613.5214
$ws.Range("N138").Value = -29139.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3446.7693
$ws.Range("I63").Value = 3256.4443
$ws.Range("J63").Value = 3875
$ws.Range("K63").Value = 3256.4443
$ws.Range("L63").Value = 3875
$ws.Range("M63").Value = -2570.4443
$ws.Range("N63").Value = -5247

$ws.Range("H66").Value = 3446.7693
$ws.Range("I66").Value = 3256.4443
$ws.Range("J66").Value = 3875
$ws.Range("K66").Value = 16282.2215
$ws.Range("L66").Value = 19375
$ws.Range("M66").Value = -12850.2215
$ws.Range("N66").Value = -26239

$ws.Range("H132").Value = 17554.791
$ws.Range("I132").Value = 22218.6
$ws.Range("K132").Value = 66655.79999999999
$ws.Range("M132").Value = -64125.79999999999

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 4500
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -2330

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3263.4285
$ws.Range("I31").Value = 2338.303
$ws.Range("J31").Value = 6655.5557
$ws.Range("K31").Value = 2338.303
$ws.Range("L31").Value = 6655.5557
$ws.Range("M31").Value = -2043.303
$ws.Range("N31").Value = -7245.5557

$ws.Range("H34").Value = 3263.4285
$ws.Range("I34").Value = 2338.303
$ws.Range("J34").Value = 6655.5557
$ws.Range("K34").Value = 2338.303
$ws.Range("L34").Value = 6655.5557
$ws.Range("M34").Value = -2136.303
$ws.Range("N34").Value = -7059.5557

$ws.Range("H107").Value = 341.10526
$ws.Range("I107").Value = 182.5
$ws.Range("J107").Value = 613
$ws.Range("K107").Value = 182.5
$ws.Range("L107").Value = 613
$ws.Range("M107").Value = 1737.5
$ws.Range("N107").Value = -4453

$ws.Range("H134").Value = 1672.0526
$ws.Range("I134").Value = 1050.9565
$ws.Range("J134").Value = 2624.4
$ws.Range("K134").Value = 3152.8695
$ws.Range("L134").Value = 7873.200000000001
$ws.Range("M134").Value = -617.8694999999998
$ws.Range("N134").Value = -12943.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 623.9091
$ws.Range("I5").Value = 623.9091
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1871.7273
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1759.7273
$ws.Range("N5").ClearContents()

$ws.Range("H38").Value = 90.34999999999999
$ws.Range("I38").Value = 65.40000000000001
$ws.Range("J38").Value = 115.3
$ws.Range("K38").Value = 196.2
$ws.Range("L38").Value = 345.9
$ws.Range("M38").Value = 150.8
$ws.Range("N38").Value = -1039.9

$ws.Range("H80").Value = 980
$ws.Range("I80").Value = 900
$ws.Range("K80").Value = 2700
$ws.Range("M80").Value = -1764

$ws.Range("H83").Value = 980
$ws.Range("I83").Value = 900
$ws.Range("K83").Value = 8100
$ws.Range("M83").Value = -3420

$ws.Range("H97").Value = 1186.9166
$ws.Range("I97").Value = 434.33334
$ws.Range("J97").Value = 1437.7778
$ws.Range("K97").Value = 1303.00002
$ws.Range("L97").Value = 4313.3334
$ws.Range("M97").Value = -807.0000199999999
$ws.Range("N97").Value = -5305.3334

$ws.Range("H98").Value = 2877.5925
$ws.Range("I98").Value = 178.42857
$ws.Range("K98").Value = 535.28571
$ws.Range("M98").Value = 962.71429

$ws.Range("H107").Value = 581.8095
$ws.Range("I107").Value = 386.75
$ws.Range("J107").Value = 627.7059
$ws.Range("K107").Value = 1160.25
$ws.Range("L107").Value = 1883.1177
$ws.Range("M107").Value = 759.75
$ws.Range("N107").Value = -5723.117700000001

$ws.Range("H117").Value = 967.1667
$ws.Range("I117").Value = 126.125
$ws.Range("J117").Value = 1640
$ws.Range("K117").Value = 378.375
$ws.Range("L117").Value = 4920
$ws.Range("M117").Value = 3063.625
$ws.Range("N117").Value = -11804

$ws.Range("H122").Value = 20834284
$ws.Range("I122").Value = 71428940
$ws.Range("J122").Value = 1188.4117
$ws.Range("K122").Value = 642860460
$ws.Range("L122").Value = 10695.7053
$ws.Range("M122").Value = -642858010
$ws.Range("N122").Value = -15595.7053

$ws.Range("H135").Value = 623.9091
$ws.Range("I135").Value = 623.9091
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5615.1819
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3080.1819
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 159195.38
$ws.Range("I122").Value = 250921.17
$ws.Range("K122").Value = 752763.51
$ws.Range("M122").Value = -750313.51

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3640
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 2733.3333
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 2733.3333
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -2957.3333

$ws.Range("H16").Value = 3466.6667
$ws.Range("I16").Value = 3560
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 3560
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -3390
$ws.Range("N16").Value = -3340

$ws.Range("H22").Value = 705.7857
$ws.Range("I22").Value = 838.5
$ws.Range("J22").Value = 606.25
$ws.Range("K22").Value = 838.5
$ws.Range("L22").Value = 606.25
$ws.Range("M22").Value = -543.5
$ws.Range("N22").Value = -1196.25

$ws.Range("H27").Value = 705.7857
$ws.Range("I27").Value = 838.5
$ws.Range("J27").Value = 606.25
$ws.Range("K27").Value = 838.5
$ws.Range("L27").Value = 606.25
$ws.Range("M27").Value = -731.5
$ws.Range("N27").Value = -820.25

$ws.Range("H46").Value = 1106.5
$ws.Range("I46").Value = 1172.5
$ws.Range("J46").Value = 842.5
$ws.Range("K46").Value = 1172.5
$ws.Range("L46").Value = 842.5
$ws.Range("M46").Value = -984.5
$ws.Range("N46").Value = -1218.5

$ws.Range("H55").Value = 208.68182
$ws.Range("I55").Value = 194.05882
$ws.Range("J55").Value = 258.4
$ws.Range("K55").Value = 194.05882
$ws.Range("L55").Value = 258.4
$ws.Range("M55").Value = -21.05882
$ws.Range("N55").Value = -604.4

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()

$ws.Range("H122").Value = 2476
$ws.Range("I122").Value = 2476
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7428
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4978
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3640
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 2733.3333
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 8199.999899999999
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -13139.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H122").Value = 142859600
$ws.Range("I122").Value = 200002270
$ws.Range("J122").Value = 2927.5
$ws.Range("K122").Value = 600006810
$ws.Range("L122").Value = 8782.5
$ws.Range("M122").Value = -600004360
$ws.Range("N122").Value = -13682.5

$ws.Range("H126").Value = 1112.5
$ws.Range("I126").Value = 1112.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3337.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -867.5
$ws.Range("N126").ClearContents()
